# Auto-generated: apply literal value updates scraped from the FFXIV Chocobo Leve
# profit-tracking workbook (scheduled-runner data refresh). Each cell below is a
# plain numeric literal (no formulas in this workbook) so the edit is just a set of
# direct Range.Value assignments per sheet; a handful of cells are cleared to blank
# ($null) because the refreshed row no longer has a value in that column.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 536.7091
$ws.Range("J17").Value = 372.7447
$ws.Range("L17").Value = 1118.2341
$ws.Range("N17").Value = -1454.2341
$ws.Range("H62").Value = 5878.222
$ws.Range("I62").Value = 4738
$ws.Range("K62").Value = 4738
$ws.Range("M62").Value = -4114
$ws.Range("H65").Value = 5878.222
$ws.Range("I65").Value = 4738
$ws.Range("K65").Value = 23690
$ws.Range("M65").Value = -20570
$ws.Range("H120").Value = 27761
$ws.Range("J120").Value = 27761
$ws.Range("L120").Value = 27761
$ws.Range("N120").Value = -37437
$ws.Range("H123").Value = 41802
$ws.Range("J123").Value = 41802
$ws.Range("L123").Value = 41802
$ws.Range("N123").Value = -51602

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2190.5
$ws.Range("I45").Value = 2190.5
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2190.5
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1813.5
$ws.Range("N45").Value = $null
$ws.Range("H80").Value = 32004.4
$ws.Range("J80").Value = 32004.4
$ws.Range("L80").Value = 32004.4
$ws.Range("N80").Value = -34000.4
$ws.Range("H83").Value = 32004.4
$ws.Range("J83").Value = 32004.4
$ws.Range("L83").Value = 96013.20000000001
$ws.Range("N83").Value = -105997.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1701.1538
$ws.Range("I105").Value = 1664.65
$ws.Range("J105").Value = 2139.2
$ws.Range("K105").Value = 1664.65
$ws.Range("L105").Value = 2139.2
$ws.Range("M105").Value = 82.34999999999991
$ws.Range("N105").Value = -5633.2
$ws.Range("H114").Value = 35342
$ws.Range("J114").Value = 35342
$ws.Range("L114").Value = 35342
$ws.Range("N114").Value = -44020

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3089.5173
$ws.Range("I31").Value = 1109.625
$ws.Range("K31").Value = 1109.625
$ws.Range("M31").Value = -814.625
$ws.Range("H34").Value = 3089.5173
$ws.Range("I34").Value = 1109.625
$ws.Range("K34").Value = 1109.625
$ws.Range("M34").Value = -907.625
$ws.Range("H87").Value = 24816.666
$ws.Range("J87").Value = 24816.666
$ws.Range("L87").Value = 24816.666
$ws.Range("N87").Value = -27188.666
$ws.Range("H90").Value = 24816.666
$ws.Range("J90").Value = 24816.666
$ws.Range("L90").Value = 74449.99800000001
$ws.Range("N90").Value = -86305.99800000001
$ws.Range("H99").Value = 12504551
$ws.Range("I99").Value = 33335416
$ws.Range("J99").Value = 6031.4
$ws.Range("K99").Value = 33335416
$ws.Range("L99").Value = 6031.4
$ws.Range("M99").Value = -33333918
$ws.Range("N99").Value = -9027.4
$ws.Range("H126").Value = 12504551
$ws.Range("I126").Value = 33335416
$ws.Range("J126").Value = 6031.4
$ws.Range("K126").Value = 100006248
$ws.Range("L126").Value = 18094.2
$ws.Range("M126").Value = -100003778
$ws.Range("N126").Value = -23034.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 557248.0600000001
$ws.Range("I5").Value = 546
$ws.Range("J5").Value = 954892.4399999999
$ws.Range("K5").Value = 1638
$ws.Range("L5").Value = 2864677.32
$ws.Range("M5").Value = -1526
$ws.Range("N5").Value = -2864901.32
$ws.Range("H22").Value = 2120.4
$ws.Range("J22").Value = 2120.4
$ws.Range("L22").Value = 6361.200000000001
$ws.Range("N22").Value = -6699.200000000001
$ws.Range("H27").Value = 2120.4
$ws.Range("J27").Value = 2120.4
$ws.Range("L27").Value = 6361.200000000001
$ws.Range("N27").Value = -6565.200000000001
$ws.Range("H113").Value = 4630237
$ws.Range("I113").Value = 637.3077
$ws.Range("K113").Value = 1911.9231
$ws.Range("M113").Value = 258.0769
$ws.Range("H122").Value = 3410.6333
$ws.Range("J122").Value = 3727.6538
$ws.Range("L122").Value = 33548.8842
$ws.Range("N122").Value = -38448.8842
$ws.Range("H135").Value = 557248.0600000001
$ws.Range("I135").Value = 546
$ws.Range("J135").Value = 954892.4399999999
$ws.Range("K135").Value = 4914
$ws.Range("L135").Value = 8594031.959999999
$ws.Range("M135").Value = -2379
$ws.Range("N135").Value = -8599101.959999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = $null
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = $null
$ws.Range("H113").Value = 1294.7273
$ws.Range("I113").Value = 1208.1666
$ws.Range("J113").Value = 1398.6
$ws.Range("K113").Value = 1208.1666
$ws.Range("L113").Value = 1398.6
$ws.Range("M113").Value = 961.8334
$ws.Range("N113").Value = -5738.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 6672.3335
$ws.Range("I4").Value = 6672.3335
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 6672.3335
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = $null
$ws.Range("M4").Value = -6559.3335
$ws.Range("H7").Value = 4222.1055
$ws.Range("J7").Value = 4553.5454
$ws.Range("L7").Value = 4553.5454
$ws.Range("N7").Value = -4777.5454
$ws.Range("H28").Value = 6672.3335
$ws.Range("I28").Value = 6672.3335
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 6672.3335
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = $null
$ws.Range("M28").Value = -6440.3335
$ws.Range("H37").Value = 6672.3335
$ws.Range("I37").Value = 6672.3335
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 6672.3335
$ws.Range("L37").Value = 0
$ws.Range("N37").Value = $null
$ws.Range("M37").Value = -6565.3335
$ws.Range("H122").Value = 4937.615
$ws.Range("I122").Value = 1752
$ws.Range("J122").Value = 5516.8184
$ws.Range("K122").Value = 5256
$ws.Range("L122").Value = 16550.4552
$ws.Range("M122").Value = -2806
$ws.Range("N122").Value = -21450.4552
$ws.Range("H126").Value = 4222.1055
$ws.Range("J126").Value = 4553.5454
$ws.Range("L126").Value = 13660.6362
$ws.Range("N126").Value = -18600.6362
$ws.Range("H132").Value = 6126.364
$ws.Range("I132").Value = 1339.8
$ws.Range("J132").Value = 7534.1763
$ws.Range("K132").Value = 4019.4
$ws.Range("L132").Value = 22602.5289
$ws.Range("M132").Value = -1489.4
$ws.Range("N132").Value = -27662.5289
$ws.Range("H136").Value = 2883.3157
$ws.Range("I136").Value = 791.8095
$ws.Range("J136").Value = 5466.9414
$ws.Range("K136").Value = 2375.4285
$ws.Range("L136").Value = 16400.8242
$ws.Range("M136").Value = 174.5715
$ws.Range("N136").Value = -21500.8242

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 39999
$ws.Range("J28").Value = 39999
$ws.Range("L28").Value = 39999
$ws.Range("N28").Value = -40695
$ws.Range("H57").Value = 19200
$ws.Range("J57").Value = 19200
$ws.Range("L57").Value = 19200
$ws.Range("N57").Value = -20708
$ws.Range("H81").Value = 4363.636
$ws.Range("I81").Value = 2600
$ws.Range("J81").Value = 5833.3335
$ws.Range("K81").Value = 5200
$ws.Range("L81").Value = 11666.667
$ws.Range("M81").Value = -4139
$ws.Range("N81").Value = -13788.667
$ws.Range("H84").Value = 4363.636
$ws.Range("I84").Value = 2600
$ws.Range("J84").Value = 5833.3335
$ws.Range("K84").Value = 26000
$ws.Range("L84").Value = 58333.335
$ws.Range("M84").Value = -20696
$ws.Range("N84").Value = -68941.33499999999
$ws.Range("H113").Value = 729.38464
$ws.Range("I113").Value = 696.5833
$ws.Range("J113").Value = 1123
$ws.Range("K113").Value = 2089.7499
$ws.Range("L113").Value = 3369
$ws.Range("M113").Value = 80.2501000000002
$ws.Range("N113").Value = -7709
$ws.Range("H126").Value = 1924.2778
$ws.Range("I126").Value = 1111.2142
$ws.Range("J126").Value = 4770
$ws.Range("K126").Value = 3333.6426
$ws.Range("L126").Value = 14310
$ws.Range("M126").Value = -863.6425999999997
$ws.Range("N126").Value = -19250
$ws.Range("H132").Value = 12349020
$ws.Range("I132").Value = 2068.3076
$ws.Range("K132").Value = 6204.9228
$ws.Range("M132").Value = -3674.9228
